# ---------------------------------------------------------------------------
# Saptamana 29 -> Saptamana 30/29/28/27 update
#   - appends weekly totals on Sheet1 (row 46 E/F, new row 47 "-")
#   - adds a right-hand commentary column (I2:I4) on Sheet1
#   - adds a brand-new "Sheet2" with the case/death characteristics table
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# xlCenter
$xlCenter = -4108

# ---------------------------------------------------------------------------
# Sheet1 updates
# ---------------------------------------------------------------------------

# header rows re-centred (B1:G1 and G2)
$ws1.Range("B1:G1").HorizontalAlignment = $xlCenter
$ws1.Range("G2").HorizontalAlignment = $xlCenter

# new commentary column I
$ws1.Columns.Item(9).ColumnWidth = 44

$ws1.Range("I2").Value = "94% din decese aveau comorbiditati asociate"
$ws1.Range("I3").Value = "78.2% din decese au fost persoane de peste 60 ani"
$ws1.Range("I4").Value = "59.7% din decese au fost barbati"
$ws1.Range("I2:I4").HorizontalAlignment = $xlCenter

# weekly totals filled in on the "Total" row
$ws1.Range("E46").Value = 7763
$ws1.Range("F46").Value = 168

# new trailing "-" row
$ws1.Range("C47").Value = "-"
$ws1.Range("D47").Value = "-"
$ws1.Range("E47").Value = "-"
$ws1.Range("F47").Value = "-"

# restore the view state on Sheet1 (selection + scroll)
$ws1.Activate()
$ws1.Range("H40").Select()

# ---------------------------------------------------------------------------
# add Sheet2 (after Sheet1) with the cases/deaths characteristics table
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# column widths
$ws2.Columns.Item(1).ColumnWidth = 13.29
$ws2.Columns.Item(2).ColumnWidth = 22.71
$ws2.Columns.Item(3).ColumnWidth = 19.86
$ws2.Columns.Item(4).ColumnWidth = 11
$ws2.Columns.Item(5).ColumnWidth = 12
$ws2.Columns.Item(6).ColumnWidth = 16.71
$ws2.Columns.Item(7).ColumnWidth = 22.86
$ws2.Columns.Item(8).ColumnWidth = 15.71
$ws2.Columns.Item(9).ColumnWidth = 23.57
$ws2.Columns.Item(11).ColumnWidth = 21.14
$ws2.Columns.Item(12).ColumnWidth = 17.71
$ws2.Columns.Item(13).ColumnWidth = 11.43
$ws2.Columns.Item(14).ColumnWidth = 18.71
$ws2.Columns.Item(15).ColumnWidth = 12

# title row
$ws2.Range("A1:O1").Merge()
$ws2.Range("A1").Value = "Caracteristicile cazurilor confirmate si a deceselor"

# section headers
$ws2.Range("B3:F3").Merge()
$ws2.Range("B3").Value = "Cazuri"
$ws2.Range("G3:P3").Merge()
$ws2.Range("G3").Value = "Decese"

# column headers (row 4)
$ws2.Range("B4").Value = "Varsta,mediana(range)"
$ws2.Range("C4").Value = "Sex,masculin"
$ws2.Range("D4").Value = "Import"
$ws2.Range("E4").Value = "Vindecati"
$ws2.Range("F4").Value = "Personal sanitar"
$ws2.Range("G4").Value = "Varsta,mediana(range)"
$ws2.Range("H4").Value = "Sex,masculin"
$ws2.Range("I4").Value = "Afectiuni cardiovasculare"
$ws2.Range("J4").Value = "Diabet"
$ws2.Range("K4").Value = "Afectiuni neurologice"
$ws2.Range("L4").Value = "Afectiuni renale"
$ws2.Range("M4").Value = "Obezitate"
$ws2.Range("N4").Value = "Afectiuni Pulmonare"
$ws2.Range("O4").Value = "Neoplasm"
$ws2.Range("P4").Value = "Altele"

# row 5 - Numar
$ws2.Range("A5").Value = "Numar"
$ws2.Range("B5").Value = "48(0-99)"
$ws2.Range("C5").Value = 21147
$ws2.Range("D5").Value = 974
$ws2.Range("E5").Value = 27754
$ws2.Range("F5").Value = 3714
$ws2.Range("G5").Value = "69(20-99)"
$ws2.Range("H5").Value = 1318
$ws2.Range("I5").Value = 1480
$ws2.Range("J5").Value = 710
$ws2.Range("K5").Value = 492
$ws2.Range("L5").Value = 445
$ws2.Range("M5").Value = 405
$ws2.Range("N5").Value = 374
$ws2.Range("O5").Value = 260
$ws2.Range("P5").Value = 436

# row 6 - Procentaj
$ws2.Range("A6").Value = "Procentaj"
$ws2.Range("C6").Value = 46.1
$ws2.Range("D6").Value = 2.1
$ws2.Range("E6").Value = 60.5
$ws2.Range("F6").Value = 8.1
$ws2.Range("H6").Value = 59.7
$ws2.Range("I6").Value = 67.1
$ws2.Range("J6").Value = 32.2
$ws2.Range("K6").Value = 22.3
$ws2.Range("L6").Value = 20.2
$ws2.Range("M6").Value = 18.4
$ws2.Range("N6").Value = 17
$ws2.Range("O6").Value = 11.8
$ws2.Range("P6").Value = 19.8

# whole table centred
$ws2.Range("A1:P6").HorizontalAlignment = $xlCenter

# view state on Sheet2
$ws2.Range("F9").Select()

# re-activate Sheet1 so it stays the selected tab, matching the source file
$ws1.Activate()
